$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Story estimate (B3) is now a rolled-up sum of the planned hours (E3:E12)
# instead of a hard-coded number.
$ws.Range("B3").Formula = "=SUM(E3:E12)"

# Task-4 (row 6) planned/burnt hours were updated.
$ws.Range("E6").Value = 17
$ws.Range("F6").Value = 8

# Scroll the view down and leave the active cell on F7, matching the
# author's last on-screen position when the sheet was saved.
$ws.Range("F7").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
